# Adds a 5th data row (TCNo 4 / URL 4 / Name "Swarup4" / Password "Swarup4")
# to the "Marta" sheet, matching the author's "adding two excel file and
# HelloAnkurSwarup1" commit: new shared string "Swarup4", a left+right thin
# border style applied to the first/last cell of the new row, and the
# worksheet selection left on F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 5)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Swarup4"
$ws.Range("D5").Value = "Swarup4"

# A5 gets a new style: thin border on the left & right edges only.
$a5 = $ws.Range("A5")
$a5.Borders.Item(7).Weight = 2
$a5.Borders.Item(10).Weight = 2

# Give D5 the exact same style as A5 by copying formats across (writing the
# same border edges again on a second, disjoint cell would otherwise mint a
# duplicate style entry instead of reusing A5's).
$a5.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the same cell selected as in the authored workbook.
[void]$ws.Range("F8").Select()

Write-Output "Added row 5 (Swarup4) with left/right borders on A5 & D5"
